$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 11
}

$newB = @{
    21 = 7
    22 = 22
    23 = 45
    24 = 81
    25 = 121
    26 = 210
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
}

foreach ($r in $newB.Keys) {
    $ws.Cells.Item($r, 2).Value = $newB[$r]
}
